$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.966.71'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '2.546.92'
$ws.Range('E3').Value = '  +3.45%  '
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '569.71'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.67%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '146.18'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +2.11%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '2.547.26'
$ws.Range('E9').Value = '  +3.47%  '
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('E11').Value = '  -3.26%  '
$ws.Range('E12').Value = '  -0.17%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.353'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.93%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '27.35'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').Value = '3.004.05'
$ws.Range('E15').Value = '  +3.60%  '
$ws.Range('D16').Value = '62.915.84'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('D18').Value = '2.543.26'
$ws.Range('E18').Value = '  +3.46%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.31'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.34'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +1.09%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '334.78'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.76%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.79'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('E23').Value = '  +0.07%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '65.18'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('E25').Value = '  -0.91%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.59'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +5.59%  '
$ws.Range('E27').Value = '  +0.01%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.47'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +2.07%  '
$ws.Range('E29').Value = '  +3.27%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '7.34'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +8.72%  '
$ws.Range('D31').Value = '0.0₃0814'
$ws.Range('E31').Value = '  +2.49%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.84'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.71%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '175.49'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -0.41%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.54'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +1.63%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '406.21'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +3.88%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.400'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.65%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '19.08'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +1.47%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('E40').Value = '  +0.57%  '
$ws.Range('E41').Value = '  +0.11%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '39.61'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -1.28%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '152.11'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +1.52%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '3.76'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.32%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '20.78'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.61%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0532'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +3.07%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.602'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.78%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0965'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.00%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.0239'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +4.11%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '18.30'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.58%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.74'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -2.92%  '
